# [TEST SCRAPE] updated files from azure vm
#
# 1. Clear the stray empty B2 cell on "ODI Batting".
# 2. Add a new "ODI Batting Extra" sheet (4th sheet) with header row
#    (copied formatting from an existing header cell) and two data rows.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the empty inline-string cell at ODI Batting!B2 ---------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").Value = ""

# --- 2. Add the new "ODI Batting Extra" worksheet at the end ---------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row
$newSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$newSheet.Cells.Item(1, 2).Value = "BATTING_POSITION"
$newSheet.Cells.Item(1, 3).Value = "NUM_4"
$newSheet.Cells.Item(1, 4).Value = "NUM_6"
$newSheet.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Data rows (MATCH_CODE is text, matching the other sheets' convention)
$newSheet.Range("A2:A3").NumberFormat = "@"
$newSheet.Cells.Item(2, 1).Value = "4421"
$newSheet.Cells.Item(2, 6).Value = "NO"

$newSheet.Cells.Item(3, 1).Value = "4460"

# Re-use the existing bold/bordered header style (style index 1 in the
# original workbook) instead of creating a fresh duplicate style.
$headerStyleSource = $wb.Worksheets.Item("Player Info").Range("A1")
$headerStyleSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the original active sheet/tab (adding a sheet would otherwise make
# the new sheet active).
$wb.Worksheets.Item("Player Info").Activate()
